$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: row 5 ("The number of children involved in PEC services (persons)")
# held decimal percentage-like values with a one-decimal number format; it
# should show whole-number counts of children, formatted as thousand-grouped
# integers (same numFmt used elsewhere in the sheet, e.g. D6:H6).
$ws.Range("D5").Value = 661
$ws.Range("D5").NumberFormat = "#\ ##0"

$ws.Range("E5").Value = 619
$ws.Range("E5").NumberFormat = "#\ ##0"

$ws.Range("F5").Value = 650
$ws.Range("F5").NumberFormat = "#\ ##0"

$ws.Range("G5").Value = 657
$ws.Range("G5").NumberFormat = "#\ ##0"

$ws.Range("H5").Value = 616
$ws.Range("H5").NumberFormat = "#\ ##0"

# Leave the selection where the editor left it after fixing the row.
$ws.Range("D5:H5").Select()
